$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '69.360.09'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.16%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.687.50'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.18%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '678.40'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.54%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '159.05'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.89%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -0.53%  '
$ws.Range('E9').Value = '  -1.28%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.15'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -3.03%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.441'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.34%  '
$ws.Range('E12').Value = '  -2.44%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.309.10'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.17%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '32.39'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.26%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.674.34'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.38%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '69.312.96'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.29%  '
$ws.Range('E17').Value = '  +2.77%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '15.99'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.67%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.48'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.89%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '468.27'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.96%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '9.87'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.89%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.653'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.99%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '80.01'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.834.15'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.11%  '
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('E26').Value = '  -5.20%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.17'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.55%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.69'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.86%  '
$ws.Range('E30').Value = '  -3.88%  '
$ws.Range('E31').Value = '  -3.41%  '
$ws.Range('B32').Value = 'Binance-PegBSC-USD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.00'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.23%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.99'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -3.33%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '26.88'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.28%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.677.70'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.55%  '
$ws.Range('E36').Value = '  -4.31%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '8.32'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.81%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '6.28'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.31%  '
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.24'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -3.49%  '
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0904'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.83%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '170.55'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +4.26%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.943'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.13%  '
$ws.Range('E45').Value = '  -1.84%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '28.08'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -7.03%  '
$ws.Range('E47').Value = '  -1.33%  '
$ws.Range('E48').Value = '  -3.36%  '
$ws.Range('B49').Value = 'FLOKI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.000276'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.33%  '
$ws.Range('B50').Value = 'ONDO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.30'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -2.56%  '
$ws.Range('E51').Value = '  -3.07%  '
